$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cupons")

# Add a new coupon row ("GOMES") right below the existing data (row 12 -> row 13)
$newRow = 13
$ws.Cells.Item($newRow, 1).Value = "GOMES"
$ws.Cells.Item($newRow, 2).Value = "porcentagem"
$ws.Cells.Item($newRow, 3).Value = 5
$ws.Cells.Item($newRow, 4).Value = Get-Date -Year 2025 -Month 12 -Day 31 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Cells.Item($newRow, 4).NumberFormat = "mm-dd-yy"
$ws.Cells.Item($newRow, 5).Value = 100
$ws.Cells.Item($newRow, 6).Value = 100
$ws.Cells.Item($newRow, 7).Value = 50
$ws.Cells.Item($newRow, 9).Value = $true
